$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Status column H2): was a numeric 0.4 (40%), now reads the "done" status text.
$ws.Cells.Item(2, 8).Value = "done"

# Row 4 grew taller (e.g. after the Notes/status text re-wrapped onto a 4th line).
$ws.Rows.Item(4).RowHeight = 58

# Rows 10, 11 and 16 move from "not started" to "done".
$ws.Cells.Item(10, 8).Value = "done"
$ws.Cells.Item(11, 8).Value = "done"
$ws.Cells.Item(16, 8).Value = "done"

# Cursor ends up parked on H1.
$ws.Range("H1").Select() | Out-Null
